$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '48.153.80'
$ws.Range('E2').Value = '  +1.87%  '

$ws.Range('D3').Value = '2.509.92'
$ws.Range('E3').Value = '  +0.70%  '

$ws.Range('E4').Value = '  +0.02%  '

Set-TextValue 'D5' '321.37'
$ws.Range('E5').Value = '  -0.12%  '

Set-TextValue 'D6' '108.60'
$ws.Range('E6').Value = '  +0.17%  '

Set-TextValue 'D7' '0.527'
$ws.Range('E7').Value = '  +0.80%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  +0.44%  '

Set-TextValue 'D10' '39.84'
$ws.Range('E10').Value = '  +1.62%  '

Set-TextValue 'D11' '20.21'
$ws.Range('E11').Value = '  +10.10%  '

$ws.Range('E12').Value = '  +0.98%  '

$ws.Range('E13').Value = '  +0.42%  '

$ws.Range('E14').Value = '  +0.82%  '

$ws.Range('D15').Value = '2.901.29'
$ws.Range('E15').Value = '  +0.70%  '

$ws.Range('D16').Value = '2.513.66'
$ws.Range('E16').Value = '  +0.94%  '

$ws.Range('E17').Value = '  -0.12%  '

$ws.Range('D18').Value = '48.003.00'
$ws.Range('E18').Value = '  +1.75%  '

Set-TextValue 'D19' '13.12'
$ws.Range('E19').Value = '  -0.33%  '

$ws.Range('E20').Value = '  +0.43%  '

$ws.Range('E21').Value = '  +1.33%  '

Set-TextValue 'D22' '2.75'
$ws.Range('E22').Value = '  +1.02%  '

Set-TextValue 'D23' '72.12'

Set-TextValue 'D24' '276.10'
$ws.Range('E24').Value = '  +12.45%  '

$ws.Range('E25').Value = '  +0.45%  '

$ws.Range('E26').Value = '  -0.09%  '

Set-TextValue 'D27' '25.89'
$ws.Range('E27').Value = '  +0.67%  '

Set-TextValue 'D28' '9.88'
$ws.Range('E28').Value = '  -0.84%  '

Set-TextValue 'D29' '35.36'
$ws.Range('E29').Value = '  +2.08%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D30' '2.10'
$ws.Range('E30').Value = '  -7.54%  '

$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D31' '0.136'
$ws.Range('E31').Value = '  -1.62%  '

Set-TextValue 'D32' '49.22'
$ws.Range('E32').Value = '  -1.12%  '

Set-TextValue 'D33' '19.41'
$ws.Range('E33').Value = '  -3.61%  '

$ws.Range('E34').Value = '  +0.32%  '

$ws.Range('E35').Value = '  -0.03%  '

$ws.Range('E36').Value = '  +0.24%  '

$ws.Range('E37').Value = '  -0.21%  '

$ws.Range('E38').Value = '  -3.11%  '

$ws.Range('E39').Value = '  +0.59%  '

Set-TextValue 'D40' '123.21'
$ws.Range('E40').Value = '  +4.22%  '

$ws.Range('E41').Value = '  +0.20%  '

$ws.Range('E42').Value = '  +0.00%  '

Set-TextValue 'D43' '21.66'
$ws.Range('E43').Value = '  -5.85%  '

$ws.Range('E44').Value = '  +3.06%  '

$ws.Range('D45').Value = '2.000.72'
$ws.Range('E45').Value = '  +0.14%  '

$ws.Range('E46').Value = '  +3.29%  '

$ws.Range('E47').Value = '  +4.55%  '

$ws.Range('E48').Value = '  -0.84%  '

Set-TextValue 'D49' '9.03'

Set-TextValue 'D50' '5.20'
$ws.Range('E50').Value = '  +2.20%  '

Set-TextValue 'D51' '79.72'
$ws.Range('E51').Value = '  +2.43%  '
